$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "25.902.12"
$ws.Range("E2").Value2 = "  +0.95%  "

$ws.Range("D3").Value2 = "1.641.86"
$ws.Range("E3").Value2 = "  +1.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.003"
$ws.Range("E4").Value2 = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "215.83"
$ws.Range("E5").Value2 = "  +0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.5081"
$ws.Range("E6").Value2 = "  +0.40%  "

$ws.Range("E7").Value2 = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.2604"
$ws.Range("E8").Value2 = "  +1.92%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06473"
$ws.Range("E9").Value2 = "  +2.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "20.34"
$ws.Range("E10").Value2 = "  +5.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.07819"
$ws.Range("E11").Value2 = "  +0.76%  "

$ws.Range("E12").Value2 = "  +1.13%  "

$ws.Range("B13").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value2 = "1.868.97"
$ws.Range("E13").Value2 = "  +1.67%  "

$ws.Range("B14").Value2 = "WrappedEther"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value2 = "1.642.67"
$ws.Range("E14").Value2 = "  +1.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.5666"
$ws.Range("E15").Value2 = "  +2.37%  "

$ws.Range("D16").Value2 = "0.0₅7716"
$ws.Range("E16").Value2 = "  +3.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "63.62"
$ws.Range("E17").Value2 = "  +0.45%  "

$ws.Range("D18").Value2 = "25.922.86"

$ws.Range("E19").Value2 = "  +0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "194.87"
$ws.Range("E20").Value2 = "  +1.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "4.409"
$ws.Range("E21").Value2 = "  +1.80%  "

$ws.Range("E22").Value2 = "  +3.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "6.286"
$ws.Range("E23").Value2 = "  +5.81%  "

$ws.Range("E24").Value2 = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "1.764"
$ws.Range("E25").Value2 = "  -3.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "139.11"
$ws.Range("E26").Value2 = "  -0.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.1231"
$ws.Range("E27").Value2 = "  -1.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "6.871"
$ws.Range("E28").Value2 = "  +2.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "15.60"
$ws.Range("E29").Value2 = "  +1.81%  "

$ws.Range("E30").Value2 = "  +1.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.05004"
$ws.Range("E31").Value2 = "  +3.60%  "

$ws.Range("E32").Value2 = "  +1.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "3.276"
$ws.Range("E33").Value2 = "  +3.60%  "

$ws.Range("E34").Value2 = "  +2.93%  "

$ws.Range("E35").Value2 = "  +0.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.9115"
$ws.Range("E36").Value2 = "  +2.74%  "

$ws.Range("E37").Value2 = "  +2.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.5549"
$ws.Range("E38").Value2 = "  +1.82%  "

$ws.Range("D39").Value2 = "1.130.84"
$ws.Range("E39").Value2 = "  +0.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.01576"
$ws.Range("E40").Value2 = "  +1.55%  "

$ws.Range("E41").Value2 = "  +0.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "5.519"
$ws.Range("E42").Value2 = "  -0.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "99.91"
$ws.Range("E43").Value2 = "  +3.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.8019"
$ws.Range("E44").Value2 = "  +1.37%  "

$ws.Range("E45").Value2 = "  -0.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "55.81"
$ws.Range("E46").Value2 = "  +2.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.4236"
$ws.Range("E47").Value2 = "  -3.84%  "

$ws.Range("B48").Value2 = "Cronos"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.05048"
$ws.Range("E48").Value2 = "  -0.84%  "

$ws.Range("B49").Value2 = "EnergySwap"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "7.682"
$ws.Range("E49").Value2 = "  +1.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.003"
$ws.Range("E50").Value2 = "  +0.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "1.002"
$ws.Range("E51").Value2 = "  +0.12%  "

